$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.178.34'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '3.909.52'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''487.69'
$ws.Range("E5").Value = '  +3.62%  '
$ws.Range("D6").Value = '''147.02'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '  -0.68%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("D11").Value = '''0.0000345'
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").Value = '''42.91'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").Value = '''10.71'
$ws.Range("E13").Value = '  +2.57%  '
$ws.Range("D14").Value = '4.528.64'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '3.912.57'
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("D16").Value = '''14.11'
$ws.Range("E16").Value = '  -7.23%  '
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -2.51%  '
$ws.Range("D20").Value = '68.256.78'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = '''430.68'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").Value = '''3.57'
$ws.Range("E22").Value = '  +4.75%  '
$ws.Range("D23").Value = '''14.93'
$ws.Range("E23").Value = '  +3.28%  '
$ws.Range("D24").Value = '''87.40'
$ws.Range("D25").Value = '''11.44'
$ws.Range("E25").Value = '  +15.34%  '
$ws.Range("D26").Value = '''11.31'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '''38.18'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").Value = '''726.20'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").Value = '''13.88'
$ws.Range("E31").Value = '  +2.67%  '
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("E33").Value = '  +3.70%  '
$ws.Range("D34").Value = '''6.31'
$ws.Range("E34").Value = '  +17.59%  '
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("D36").Value = '0.0₃0873'
$ws.Range("E36").Value = '  +3.36%  '
$ws.Range("D37").Value = '''60.27'
$ws.Range("E37").Value = '  +3.90%  '
$ws.Range("D38").Value = '''0.409'
$ws.Range("E38").Value = '  +21.48%  '
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("D40").Value = '''0.997'
$ws.Range("E40").Value = '  -0.27%  '
$ws.Range("D41").Value = '''2.97'
$ws.Range("E41").Value = '  +16.85%  '
$ws.Range("D42").Value = '''0.0480'
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").Value = '''3.17'
$ws.Range("E43").Value = '  +3.64%  '
$ws.Range("E44").Value = '  +3.16%  '
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.33'
$ws.Range("E46").Value = '  +4.75%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '''1.00'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0342'
$ws.Range("E50").Value = '  +30.51%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '''143.93'
$ws.Range("E51").Value = '  -2.48%  '
